{"js": "// Replace the date line and every \"a\u00d7b=\" expression in the multiplication\n// table with the new values from the commit.\nconst replacements = [\n  [\"2024-01-29 Monday\", \"2024-01-30 Tuesday\"],\n  [\"66\u00d754=\", \"26\u00d733=\"],\n  [\"70\u00d752=\", \"20\u00d794=\"],\n  [\"48\u00d768=\", \"71\u00d720=\"],\n  [\"60\u00d757=\", \"86\u00d761=\"],\n  [\"88\u00d794=\", \"63\u00d721=\"],\n  [\"67\u00d771=\", \"37\u00d772=\"],\n  [\"75\u00d740=\", \"73\u00d756=\"],\n  [\"67\u00d789=\", \"95\u00d760=\"],\n  [\"31\u00d764=\", \"22\u00d770=\"],\n  [\"53\u00d731=\", \"72\u00d730=\"],\n  [\"80\u00d762=\", \"99\u00d712=\"],\n  [\"37\u00d732=\", \"21\u00d740=\"],\n  [\"50\u00d718=\", \"93\u00d770=\"],\n  [\"13\u00d730=\", \"96\u00d792=\"],\n  [\"30\u00d717=\", \"42\u00d725=\"],\n  [\"54\u00d719=\", \"43\u00d779=\"],\n  [\"74\u00d724=\", \"32\u00d740=\"],\n  [\"44\u00d769=\", \"30\u00d746=\"],\n  [\"56\u00d739=\", \"48\u00d772=\"],\n  [\"24\u00d724=\", \"15\u00d769=\"],\n  [\"19\u00d724=\", \"68\u00d748=\"],\n  [\"20\u00d799=\", \"49\u00d729=\"],\n  [\"20\u00d774=\", \"33\u00d737=\"],\n  [\"32\u00d722=\", \"62\u00d735=\"],\n  [\"30\u00d771=\", \"40\u00d779=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"a\u00d7b=\" expression in the multiplication\n# table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-01-29 Monday\", \"2024-01-30 Tuesday\"),\n    @(\"66\u00d754=\", \"26\u00d733=\"),\n    @(\"70\u00d752=\", \"20\u00d794=\"),\n    @(\"48\u00d768=\", \"71\u00d720=\"),\n    @(\"60\u00d757=\", \"86\u00d761=\"),\n    @(\"88\u00d794=\", \"63\u00d721=\"),\n    @(\"67\u00d771=\", \"37\u00d772=\"),\n    @(\"75\u00d740=\", \"73\u00d756=\"),\n    @(\"67\u00d789=\", \"95\u00d760=\"),\n    @(\"31\u00d764=\", \"22\u00d770=\"),\n    @(\"53\u00d731=\", \"72\u00d730=\"),\n    @(\"80\u00d762=\", \"99\u00d712=\"),\n    @(\"37\u00d732=\", \"21\u00d740=\"),\n    @(\"50\u00d718=\", \"93\u00d770=\"),\n    @(\"13\u00d730=\", \"96\u00d792=\"),\n    @(\"30\u00d717=\", \"42\u00d725=\"),\n    @(\"54\u00d719=\", \"43\u00d779=\"),\n    @(\"74\u00d724=\", \"32\u00d740=\"),\n    @(\"44\u00d769=\", \"30\u00d746=\"),\n    @(\"56\u00d739=\", \"48\u00d772=\"),\n    @(\"24\u00d724=\", \"15\u00d769=\"),\n    @(\"19\u00d724=\", \"68\u00d748=\"),\n    @(\"20\u00d799=\", \"49\u00d729=\"),\n    @(\"20\u00d774=\", \"33\u00d737=\"),\n    @(\"32\u00d722=\", \"62\u00d735=\"),\n    @(\"30\u00d771=\", \"40\u00d779=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
